$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new logbook entries for Week 4 and Week 5.
# Values are set in the same order the author originally typed them so that
# the shared-string table indices line up with the target workbook
# (rows 39,38,37,36 then 44 then 46,47,48,49).

$ws.Range("B39").Value2 = "Brief meeting with the hardware-crew regarding what we had done during the week. Reviewed Eddies code and showed everyone how we should work with branches and pullrequests."
$ws.Range("B39").WrapText = $true
$ws.Rows.Item(39).RowHeight = 45

$ws.Range("B38").Value2 = "Discussed positioning with Kiro. Explained how my solution worked."

$ws.Range("B37").Value2 = "Short talk with Kiro reviewing and discussing the code he had written. Followed by a Stand-up meeting with the entire group. "
$ws.Range("B37").WrapText = $true
$ws.Rows.Item(37).RowHeight = 30

$ws.Range("B36").Value2 = "Meeting with the hardware-crew plus Dejan. Discussed how the communication between units (rpi, arduino, app) will work."
$ws.Range("B36").WrapText = $true
$ws.Rows.Item(36).RowHeight = 30

$ws.Range("B44").Value2 = "Sprint-planning with the project group."

$ws.Range("B46").Value2 = "Met with Eddie at school and tried out the bluetooth connection to the application. Did not work. Looked for other solutions, landed on a Bluetooth classic solution where a server socket is created on the rpi. Had a chat with Micke from the frontend, we will meet up tomorrow to further develop our bluetooth solution."
$ws.Range("B46").WrapText = $true
$ws.Rows.Item(46).RowHeight = 75

$ws.Range("B47").Value2 = "Met with Eddie and Micke in school. Had some tech-trouble with the rpi so installed a new raspianOS on a separate SD-card. Confirmed connection between their app and the rpi with Bluetooth classic."
$ws.Range("B47").WrapText = $true
$ws.Rows.Item(47).RowHeight = 45

$ws.Range("B48").Value2 = "Meeting with the entire group. Discussed what we should be doing going forward. Gave the Rpi to Kiro as he and a member of the front-end team wanted to work with it over the weekend."
$ws.Range("B48").WrapText = $true
$ws.Rows.Item(48).RowHeight = 45

$ws.Range("B49").Value2 = "I finished writing the code for the rpi-controls. It should now be able to receive messages from a bluetooth client socket and through that control the Robots state-machine as well as listening to the robots actions during autonomous driving. Also calculating the position when moving forward. Need to do some testing to make sure that it works as intended during next week."
$ws.Range("B49").WrapText = $true
$ws.Rows.Item(49).RowHeight = 90

# Update the saved view/selection to match where the author left off scrolling.
$win = $excel.ActiveWindow
$ws.Range("C47").Select()
$win.ScrollRow = 31
$win.ScrollColumn = 1
